$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.59486203616575
$ws.Range("C2").Value = 8.791737470398171
$ws.Range("E2").Value = 15.44322561161456
$ws.Range("F2").Value = 39.49423673964566
$ws.Range("G2").Value = 3.671879754973645
$ws.Range("I2").Value = 26.51250918074938
$ws.Range("J2").Value = 8.411503718384983
$ws.Range("K2").Value = 9.376781223471033
$ws.Range("L2").Value = 12.47453874763217
$ws.Range("M2").Value = 15.80215168760568
$ws.Range("N2").Value = 20.71432252729305
$ws.Range("O2").Value = 26.30825110385946
$ws.Range("B3").Value = 12.39988198576366
$ws.Range("C3").Value = 8.778431492385291
$ws.Range("E3").Value = 15.45705149422015
$ws.Range("F3").Value = 39.53619093940004
$ws.Range("G3").Value = 3.673453175666221
$ws.Range("I3").Value = 26.59887815746057
$ws.Range("J3").Value = 8.399164393675292
$ws.Range("K3").Value = 9.240595779168917
$ws.Range("L3").Value = 12.47319787936689
$ws.Range("M3").Value = 15.76839360776916
$ws.Range("N3").Value = 20.77152126345495
$ws.Range("O3").Value = 26.38904526324436
$ws.Range("B4").Value = 12.28070698149535
$ws.Range("C4").Value = 8.770189220713057
$ws.Range("E4").Value = 15.46744236290957
$ws.Range("F4").Value = 39.56976264717586
$ws.Range("G4").Value = 3.674471143303494
$ws.Range("I4").Value = 26.65586542160723
$ws.Range("J4").Value = 8.391452877780441
$ws.Range("K4").Value = 9.157366939139317
$ws.Range("L4").Value = 12.47393133733906
$ws.Range("M4").Value = 15.74969913536738
$ws.Range("N4").Value = 20.80830598352559
$ws.Range("O4").Value = 26.4430081758367
$ws.Range("B5").Value = 12.23233889675608
$ws.Range("C5").Value = 8.766812664293564
$ws.Range("E5").Value = 15.47215578014812
$ws.Range("F5").Value = 39.58540761999171
$ws.Range("G5").Value = 3.674899058598757
$ws.Range("I5").Value = 26.68008336870131
$ws.Range("J5").Value = 8.388276644012826
$ws.Range("K5").Value = 9.123588141716125
$ws.Range("L5").Value = 12.47462308215333
$ws.Range("M5").Value = 15.7425976908972
$ws.Range("N5").Value = 20.82371579643054
$ws.Range("O5").Value = 26.46609313641599
$ws.Range("B6").Value = 12.22432105139157
$ws.Range("C6").Value = 8.766250929234069
$ws.Range("E6").Value = 15.47296739661648
$ws.Range("F6").Value = 39.58812407612128
$ws.Range("G6").Value = 3.674970905079555
$ws.Range("I6").Value = 26.68416485593196
$ws.Range("J6").Value = 8.387747202891788
$ws.Range("K6").Value = 9.117988670942031
$ws.Range("L6").Value = 12.47476171178868
$ws.Range("M6").Value = 15.74144986183828
$ws.Range("N6").Value = 20.82629997069042
$ws.Range("O6").Value = 26.46999247612525
$ws.Range("B7").Value = 12.28005379670746
$ws.Range("C7").Value = 8.770143754649434
$ws.Range("E7").Value = 15.4675039889542
$ws.Range("F7").Value = 39.56996568821745
$ws.Range("G7").Value = 3.674476861284688
$ws.Range("I7").Value = 26.65618800346414
$ws.Range("J7").Value = 8.391410178222644
$ws.Range("K7").Value = 9.156910776755582
$ws.Range("L7").Value = 12.47393907403785
$ws.Range("M7").Value = 15.74960126373198
$ws.Range("N7").Value = 20.80851210454794
$ws.Range("O7").Value = 26.44331507617258
$ws.Range("B8").Value = 12.52755601140119
$ws.Range("C8").Value = 8.787164840463999
$ws.Range("E8").Value = 15.44759860365222
$ws.Range("F8").Value = 39.50708136656692
$ws.Range("G8").Value = 3.672411526241834
$ws.Range("I8").Value = 26.54146825074519
$ws.Range("J8").Value = 8.407277311577465
$ws.Range("K8").Value = 9.329767387541969
$ws.Range("L8").Value = 12.47375418090133
$ws.Range("M8").Value = 15.79009330075104
$ws.Range("N8").Value = 20.73369991594145
$ws.Range("O8").Value = 26.33520478827657
$ws.Range("B9").Value = 13.0146075300926
$ws.Range("C9").Value = 8.819951237460836
$ws.Range("E9").Value = 15.42361438263897
$ws.Range("F9").Value = 39.44572930583825
$ws.Range("G9").Value = 3.668771267560379
$ws.Range("I9").Value = 26.34788088946775
$ws.Range("J9").Value = 8.437310465454249
$ws.Range("K9").Value = 9.670098411384476
$ws.Range("L9").Value = 12.48567891318898
$ws.Range("M9").Value = 15.8853770148377
$ws.Range("N9").Value = 20.60014465494781
$ws.Range("O9").Value = 26.15777550560068
$ws.Range("B10").Value = 13.37004796302259
$ws.Range("C10").Value = 8.843656027187603
$ws.Range("E10").Value = 15.41511295586941
$ws.Range("F10").Value = 39.43837204111758
$ws.Range("G10").Value = 3.666344115563727
$ws.Range("I10").Value = 26.22476116293986
$ws.Range("J10").Value = 8.458703654075697
$ws.Range("K10").Value = 9.918685202119315
$ws.Range("L10").Value = 12.50183839970724
$ws.Range("M10").Value = 15.96470434508189
$ws.Range("N10").Value = 20.50995984663399
$ws.Range("O10").Value = 26.04852219944356
$ws.Range("B11").Value = 13.5305273575522
$ws.Range("C11").Value = 8.854351895055736
$ws.Range("E11").Value = 15.41321240795392
$ws.Range("F11").Value = 39.44318939214641
$ws.Range("G11").Value = 3.665293115917471
$ws.Range("I11").Value = 26.17289622456474
$ws.Range("J11").Value = 8.468287351517096
$ws.Range("K11").Value = 10.03099005420769
$ws.Range("L11").Value = 12.51077413538855
$ws.Range("M11").Value = 16.002733726218
$ws.Range("N11").Value = 20.47063928914281
$ws.Range("O11").Value = 26.00340787282063
$ws.Range("B12").Value = 13.59106660909263
$ws.Range("C12").Value = 8.858389095066963
$ws.Range("E12").Value = 15.41277425668362
$ws.Range("F12").Value = 39.44618413338918
$ws.Range("G12").Value = 3.664902728518645
$ws.Range("I12").Value = 26.15385185595088
$ws.Range("J12").Value = 8.471894871104483
$ws.Range("K12").Value = 10.07336747104397
$ws.Range("L12").Value = 12.51438363072038
$ws.Range("M12").Value = 16.01740630275792
$ws.Range("N12").Value = 20.45599354680369
$ws.Range("O12").Value = 25.98698409791428
$ws.Range("B13").Value = 13.57803953859437
$ws.Range("C13").Value = 8.857520202170763
$ws.Range("E13").Value = 15.41285611989617
$ws.Range("F13").Value = 39.44548716665766
$ws.Range("G13").Value = 3.664986467835449
$ws.Range("I13").Value = 26.15792690452897
$ws.Range("J13").Value = 8.471118893724645
$ws.Range("K13").Value = 10.0642479979417
$ws.Range("L13").Value = 12.51359625772584
$ws.Range("M13").Value = 16.01423434415453
$ws.Range("N13").Value = 20.45913693003676
$ws.Range("O13").Value = 25.99049188626496
$ws.Range("B14").Value = 13.53551289947431
$ws.Range("C14").Value = 8.854684303225735
$ws.Range("E14").Value = 15.4131707260355
$ws.Range("F14").Value = 39.44341233503155
$ws.Range("G14").Value = 3.665260846332249
$ws.Range("I14").Value = 26.17131749285855
$ws.Range("J14").Value = 8.46858457939762
$ws.Range("K14").Value = 10.03447968883079
$ws.Range("L14").Value = 12.51106657963794
$ws.Range("M14").Value = 16.00393544843061
$ws.Range("N14").Value = 20.46942948909456
$ws.Range("O14").Value = 26.0020434471848
$ws.Range("B15").Value = 13.50943240882963
$ws.Range("C15").Value = 8.852945513159755
$ws.Range("E15").Value = 15.41340005742521
$ws.Range("F15").Value = 39.44229376032327
$ws.Range("G15").Value = 3.665429900208982
$ws.Range("I15").Value = 26.17959720095266
$ws.Range("J15").Value = 8.467029413974862
$ws.Range("K15").Value = 10.01622510104998
$ws.Range("L15").Value = 12.5095464089041
$ws.Range("M15").Value = 15.99766222329599
$ws.Range("N15").Value = 20.47576573601574
$ws.Range("O15").Value = 26.00920508725009
$ws.Range("B16").Value = 13.3595309831129
$ws.Range("C16").Value = 8.842955222336316
$ws.Range("E16").Value = 15.41527663086921
$ws.Range("F16").Value = 39.43822121917125
$ws.Range("G16").Value = 3.666413866784889
$ws.Range("I16").Value = 26.22823402489318
$ws.Range("J16").Value = 8.458074328781565
$ws.Range("K16").Value = 9.911326866813106
$ws.Range("L16").Value = 12.50128613561492
$ws.Range("M16").Value = 15.96225747997507
$ws.Range("N16").Value = 20.51256375386275
$ws.Range("O16").Value = 26.05156285566357
$ws.Range("B17").Value = 13.26721992728198
$ws.Range("C17").Value = 8.836803720664212
$ws.Range("E17").Value = 15.41693070994402
$ws.Range("F17").Value = 39.43781164694933
$ws.Range("G17").Value = 3.667031079156156
$ws.Range("I17").Value = 26.25913223502824
$ws.Range("J17").Value = 8.452542503548516
$ws.Range("K17").Value = 9.846748303079625
$ws.Range("L17").Value = 12.49662311577221
$ws.Range("M17").Value = 15.94103006142876
$ws.Range("N17").Value = 20.5355740617276
$ws.Range("O17").Value = 26.07872291327914
$ws.Range("B18").Value = 13.21401419157128
$ws.Range("C18").Value = 8.833257328405473
$ws.Range("E18").Value = 15.4180672984233
$ws.Range("F18").Value = 39.43834478067338
$ws.Range("G18").Value = 3.667391085787498
$ws.Range("I18").Value = 26.27729401023025
$ws.Range("J18").Value = 8.449346816348816
$ws.Range("K18").Value = 9.809533336839426
$ws.Range("L18").Value = 12.49409036692504
$ws.Range("M18").Value = 15.92900397588796
$ws.Range("N18").Value = 20.54896952491568
$ws.Range("O18").Value = 26.0947762371737
$ws.Range("B19").Value = 13.19598236176557
$ws.Range("C19").Value = 8.832055185853628
$ws.Range("E19").Value = 15.41848397453173
$ws.Range("F19").Value = 39.43865741138006
$ws.Range("G19").Value = 3.667513838111108
$ws.Range("I19").Value = 26.28351024765969
$ws.Range("J19").Value = 8.448262428909207
$ws.Range("K19").Value = 9.796921961294135
$ws.Range("L19").Value = 12.49325852702969
$ws.Range("M19").Value = 15.92496387333034
$ws.Range("N19").Value = 20.5535326048087
$ws.Range("O19").Value = 26.10028572042313
$ws.Range("B20").Value = 13.27705848654483
$ws.Range("C20").Value = 8.83745941197267
$ws.Range("E20").Value = 15.41673547126247
$ws.Range("F20").Value = 39.43777571107533
$ws.Range("G20").Value = 3.666964858376667
$ws.Range("I20").Value = 26.25580271014821
$ws.Range("J20").Value = 8.453132819488033
$ws.Range("K20").Value = 9.853630438327347
$ws.Range("L20").Value = 12.49710406571852
$ws.Range("M20").Value = 15.9432708338649
$ws.Range("N20").Value = 20.533107967403
$ws.Range("O20").Value = 26.07578700738561
$ws.Range("B21").Value = 13.54801070728404
$ws.Range("C21").Value = 8.855517634854557
$ws.Range("E21").Value = 15.41307068815615
$ws.Range("F21").Value = 39.4439900273677
$ws.Range("G21").Value = 3.66518004864767
$ws.Range("I21").Value = 26.16736818354865
$ws.Range("J21").Value = 8.469329558180352
$ws.Range("K21").Value = 10.0432277298514
$ws.Range("L21").Value = 12.51180349905957
$ws.Range("M21").Value = 16.00695317357903
$ws.Range("N21").Value = 20.46639969870872
$ws.Range("O21").Value = 25.99863255541802
$ws.Range("B22").Value = 13.72372374881666
$ws.Range("C22").Value = 8.867243175070623
$ws.Range("E22").Value = 15.4123159449791
$ws.Range("F22").Value = 39.45487196281901
$ws.Range("G22").Value = 3.664057874110104
$ws.Range("I22").Value = 26.11304365848839
$ws.Range("J22").Value = 8.479788905218671
$ws.Range("K22").Value = 10.16624961003917
$ws.Range("L22").Value = 12.52272504449871
$ws.Range("M22").Value = 16.05015300787657
$ws.Range("N22").Value = 20.424224297054
$ws.Range("O22").Value = 25.95205516296095
$ws.Range("B23").Value = 13.63008579984555
$ws.Range("C23").Value = 8.860992205193114
$ws.Range("E23").Value = 15.41256911341194
$ws.Range("F23").Value = 39.44844133293829
$ws.Range("G23").Value = 3.664652757911139
$ws.Range("I23").Value = 26.14171993146451
$ws.Range("J23").Value = 8.474218202570903
$ws.Range("K23").Value = 10.10068423710429
$ws.Range("L23").Value = 12.51677646687907
$ws.Range("M23").Value = 16.02695452354159
$ws.Range("N23").Value = 20.44660431288719
$ws.Range("O23").Value = 25.97656214286222
$ws.Range("B24").Value = 13.27261088969356
$ws.Range("C24").Value = 8.837163004466902
$ws.Range("E24").Value = 15.41682316034603
$ws.Range("F24").Value = 39.4377895633061
$ws.Range("G24").Value = 3.666994780705998
$ws.Range("I24").Value = 26.25730674865558
$ws.Range("J24").Value = 8.452865985493311
$ws.Range("K24").Value = 9.850519295821996
$ws.Range("L24").Value = 12.49688616682657
$ws.Range("M24").Value = 15.94225722611145
$ws.Range("N24").Value = 20.53422236982654
$ws.Range("O24").Value = 26.07711296413425
$ws.Range("B25").Value = 12.88302464479912
$ws.Range("C25").Value = 8.811147320907047
$ws.Range("E25").Value = 15.42849689910988
$ws.Range("F25").Value = 39.45569615393439
$ws.Range("G25").Value = 3.669712436976762
$ws.Range("I25").Value = 26.39689429396489
$ws.Range("J25").Value = 8.429302823578499
$ws.Range("K25").Value = 9.578122889949082
$ws.Range("L25").Value = 12.48114653969677
$ws.Range("M25").Value = 15.85793588928405
$ws.Range("N25").Value = 20.63487517904245
$ws.Range("O25").Value = 26.20207055731147
